$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "testpoints"

# Update lat/lon/site data (rows 2-11), columns A-D
$data = @(
    @(2, 40.81417, -96.69963, 1, "Example Site 1"),
    @(3, 33.74774, -116.927828, 2, "Example Site 2"),
    @(4, 37.2335, -112.8752, 3, "Example Site 3"),
    @(5, 30.2775, -97.82528000000001, 4, "Example Site 4"),
    @(6, 31.346534, -92.401505, 5, "Example Site 5"),
    @(7, 39.45225, -80.13771, 6, "Example Site 6"),
    @(8, 33.437222, -86.7375, 7, "Example Site 7"),
    @(9, 29.564069, -98.191041, 8, "Example Site 8"),
    @(10, 38.603834, -121.46303, 9, "Example Site 9"),
    @(11, 39.27679, -84.41392999999999, 10, "Example Site 10")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
